$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()
$ws.Cells.Item(34, 1).Value = "test"
